# DP: main - cleanup after merge
# Update cached pop/aprt_20/pop_delta/aprt_20_delta values (columns D, F, H, I)
# for the rows affected by the upstream recalculation after the merge.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @{
    21  = @{ D = 39951.38161184122;  F = 9693;               H = -4386.618388158779;  I = 0 }
    22  = @{ D = 30553.66036706011;  F = 6789.702303791135;  H = 6669.660367060107;   I = 2679.702303791135 }
    23  = @{ D = 29005.76516340315;  F = 5273.775484255118;  H = 6331.765163403154;   I = 5121.775484255118 }
    24  = @{ D = 15440.57446953674;  F = 2807.377176279408;  H = 3370.574469536741;   I = 1146.377176279408 }
    25  = @{ D = 18151.03593280999;  F = 4537.758983202499;  H = -1992.964067190005;  I = 1919.758983202499 }
    26  = @{ D = 27115.72301062605;  F = 6711;                H = -2977.276989373946;  I = 0 }
    35  = @{ D = 17799.62092021091;  F = 4034;                H = -1954.379079789087;  I = 0 }
    36  = @{ D = 4380.072759600347;  F = 1095.018189900087;  H = -480.9272403996529;  I = -353.9818100999132 }
    37  = @{ D = 13357.37473530458;  F = 3553;                H = -1466.625264695424;  I = 0 }
    350 = @{ D = 12431.98186879366;  F = 2762.662637509703;  H = -1365.018131206338;  I = 537.6626375097026 }
    351 = @{ D = 5704.637037858424;  F = 1166;                H = -626.3629621415757;  I = 0 }
    352 = @{ D = 19691.85560343674;  F = 4375.967911874831;  H = -2162.144396563261;  I = 913.9679118748309 }
    354 = @{ D = 1950.803851992338;  F = 433.5119671094083;  H = -214.1961480076625;  I = -213.4880328905917 }
    355 = @{ D = 3489.821394349341;  F = 775.5158654109646;  H = -383.1786056506594;  I = -498.4841345890354 }
    356 = @{ D = 17640.13256833902;  F = 4410.033142084756;  H = -1936.867431660976;  I = 1388.033142084756 }
    357 = @{ D = 11256.99423692391;  F = 2501.55427487198;   H = -1236.005763076088;  I = 1508.55427487198 }
    358 = @{ D = 19174.64479002168;  F = 4671;                H = -2105.355209978319;  I = 0 }
    359 = @{ D = 23556.51967789177;  F = 4283.003577798503;  H = -2586.480322108233;  I = -1815.996422201497 }
}

foreach ($rowNum in $rowUpdates.Keys) {
    $vals = $rowUpdates[$rowNum]
    $ws.Range("D$rowNum").Value = $vals.D
    $ws.Range("F$rowNum").Value = $vals.F
    $ws.Range("H$rowNum").Value = $vals.H
    $ws.Range("I$rowNum").Value = $vals.I
}
